$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch workbook calculation from manual to automatic
$excel.Calculation = -4105

# New rows 29-35 (Thoracic Surgery sessional mapping additions)

# Row 29
$ws.Cells.Item(29,1).Value = "Thoracic Surgery"
$ws.Cells.Item(29,2).Value = "1470 MADISON CANCER CENTER"
$ws.Cells.Item(29,3).Value = "8005002"
$ws.Cells.Item(29,4).Value = "LASKEY, DANIEL HENDRIK"
$ws.Cells.Item(29,5).Value = 1700206075
$ws.Cells.Item(29,6).Value = "Department"
$ws.Cells.Item(29,7).Value = "MSH- AMBULATORY CARE"
$ws.Cells.Item(29,8).Value = "Office"
$ws.Range("B29:E29").Font.Color = 3355443

# Row 30
$ws.Cells.Item(30,1).Value = "Thoracic Surgery"
$ws.Cells.Item(30,2).Value = "1470 MADISON CANCER CENTER"
$ws.Cells.Item(30,3).Value = "8005002"
$ws.Cells.Item(30,4).Value = "HOUSMAN, BRIAN NEUMANN"
$ws.Cells.Item(30,5).Value = 1649512955
$ws.Cells.Item(30,6).Value = "Department"
$ws.Cells.Item(30,7).Value = "MSH- AMBULATORY CARE"
$ws.Cells.Item(30,8).Value = "Office"
$ws.Range("B30:E30").Font.Color = 3355443

# Row 31 (no H value)
$ws.Cells.Item(31,1).Value = "Thoracic Surgery"
$ws.Cells.Item(31,2).Value = "10 E 102 PULMONARY"
$ws.Cells.Item(31,3).Value = 8005003
$ws.Cells.Item(31,4).Value = "LASKEY, DANIEL HENDRIK"
$ws.Cells.Item(31,5).Value = 1700206075
$ws.Cells.Item(31,6).Value = "Department"
$ws.Cells.Item(31,7).Value = "MSH- AMBULATORY CARE"
$ws.Range("D31:E31").Font.Color = 3355443

# Row 32 (no H value)
$ws.Cells.Item(32,1).Value = "Thoracic Surgery"
$ws.Cells.Item(32,2).Value = "10 E 102 PULMONARY"
$ws.Cells.Item(32,3).Value = 8005003
$ws.Cells.Item(32,4).Value = "SCHEININ, SCOTT ALAN"
$ws.Cells.Item(32,5).Value = 1619977238
$ws.Cells.Item(32,6).Value = "Department"
$ws.Cells.Item(32,7).Value = "MSH- AMBULATORY CARE"

# Row 33
$ws.Cells.Item(33,1).Value = "Thoracic Surgery"
$ws.Cells.Item(33,2).Value = "10 UNION SQ E RESP"
$ws.Cells.Item(33,3).Value = 8806002
$ws.Cells.Item(33,4).Value = "KAUFMAN, ANDREW J."
$ws.Cells.Item(33,5).Value = 1396907515
$ws.Cells.Item(33,6).Value = "Department"
$ws.Cells.Item(33,7).Value = "MSUS"
$ws.Cells.Item(33,8).Value = "Office"

# Row 34
$ws.Cells.Item(34,1).Value = "Thoracic Surgery"
$ws.Cells.Item(34,2).Value = "5 CUBA HILL SURGERY"
$ws.Cells.Item(34,3).Value = 8849025
$ws.Cells.Item(34,4).Value = "HOUSMAN, BRIAN NEUMANN"
$ws.Cells.Item(34,5).Value = 1649512955
$ws.Cells.Item(34,6).Value = "Department"
$ws.Cells.Item(34,7).Value = "Network"
$ws.Cells.Item(34,8).Value = "Office"

# Row 35
$ws.Cells.Item(35,1).Value = "Thoracic Surgery"
$ws.Cells.Item(35,2).Value = "30-14 CRESCENT STREET HEM ONC"
$ws.Cells.Item(35,3).Value = 8986015
$ws.Cells.Item(35,4).Value = "LEE, DONG-SEOK D."
$ws.Cells.Item(35,5).Value = 1730349085
$ws.Cells.Item(35,6).Value = "Department"
$ws.Cells.Item(35,7).Value = "Network"
$ws.Cells.Item(35,8).Value = "Office"

# Update selection to match the target workbook state
$ws.Range("B5").Select()
